$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Hyperlinks.Delete()
$ws.Range("A2:D13").ClearContents()
$ws.Range("A2:D13").Select()
